$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap row 10 and row 11 (F:V) ---
$ws.Range("F10").Value = "Floriana"
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = "Birkirkara"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 2.79
$ws.Range("K10").Value = "22/09/2023 06:13"
$ws.Range("L10").Value = 2.93
$ws.Range("M10").Value = "23/09/2023 17:58"
$ws.Range("N10").Value = 2.99
$ws.Range("O10").Value = "22/09/2023 06:13"
$ws.Range("P10").Value = 2.68
$ws.Range("Q10").Value = "23/09/2023 17:58"
$ws.Range("R10").Value = 2.37
$ws.Range("S10").Value = "22/09/2023 06:13"
$ws.Range("T10").Value = 2.8
$ws.Range("U10").Value = "23/09/2023 17:58"
$ws.Range("V10").Value = "https://www.betexplorer.com/football/malta/premier-league/floriana-birkirkara/xOKjpwxf/"

$ws.Range("F11").Value = "Hibernians"
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = "Santa Lucia"
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 1.28
$ws.Range("K11").Value = "22/09/2023 06:13"
$ws.Range("L11").Value = 1.35
$ws.Range("M11").Value = "23/09/2023 17:52"
$ws.Range("N11").Value = 5
$ws.Range("O11").Value = "22/09/2023 06:13"
$ws.Range("P11").Value = 4.84
$ws.Range("Q11").Value = "23/09/2023 17:52"
$ws.Range("R11").Value = 7.38
$ws.Range("S11").Value = "22/09/2023 06:13"
$ws.Range("T11").Value = 8.12
$ws.Range("U11").Value = "23/09/2023 17:52"
$ws.Range("V11").Value = "https://www.betexplorer.com/football/malta/premier-league/hibernians-santa-lucia/WpKfqci0/"

# --- Swap row 19 and row 20 (F:V) ---
$ws.Range("F19").Value = "Hamrun"
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = "Gudja"
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1.23
$ws.Range("K19").Value = "26/09/2023 15:13"
$ws.Range("L19").Value = 1.37
$ws.Range("M19").Value = "27/09/2023 17:05"
$ws.Range("N19").Value = 5.16
$ws.Range("O19").Value = "26/09/2023 15:13"
$ws.Range("P19").Value = 4.5
$ws.Range("Q19").Value = "27/09/2023 17:05"
$ws.Range("R19").Value = 9.19
$ws.Range("S19").Value = "26/09/2023 15:13"
$ws.Range("T19").Value = 8.43
$ws.Range("U19").Value = "27/09/2023 17:05"
$ws.Range("V19").Value = "https://www.betexplorer.com/football/malta/premier-league/hamrun-gudja/CIfTcd6g/"

$ws.Range("F20").Value = "Marsaxlokk"
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = "Sirens"
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 2.9
$ws.Range("K20").Value = "26/09/2023 12:44"
$ws.Range("L20").Value = 1.67
$ws.Range("M20").Value = "27/09/2023 16:20"
$ws.Range("N20").Value = 3.1
$ws.Range("O20").Value = "26/09/2023 12:44"
$ws.Range("P20").Value = 3.87
$ws.Range("Q20").Value = "27/09/2023 16:05"
$ws.Range("R20").Value = 2.29
$ws.Range("S20").Value = "26/09/2023 12:44"
$ws.Range("T20").Value = 4.65
$ws.Range("U20").Value = "27/09/2023 16:20"
$ws.Range("V20").Value = "https://www.betexplorer.com/football/malta/premier-league/marsaxlokk-sirens/j9gXdGLa/"

# --- Append new rows 36-42, first copying style from row 35 ---
for ($r = 36; $r -le 42; $r++) {
    $ws.Range("A35:V35").Copy($ws.Range("A" + $r + ":V" + $r))
}

# Row 36
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "malta"
$ws.Range("C36").Value = "premier-league"
$ws.Range("D36").Value = "2023-2024"
$ws.Range("E36").Value = 45227.625
$ws.Range("F36").Value = "Valletta"
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = "Hamrun"
$ws.Range("I36").Value = 2
$ws.Range("J36").Value = 4.37
$ws.Range("K36").Value = "27/10/2023 17:43"
$ws.Range("L36").Value = 4.96
$ws.Range("M36").Value = "28/10/2023 13:26"
$ws.Range("N36").Value = 3.39
$ws.Range("O36").Value = "27/10/2023 17:43"
$ws.Range("P36").Value = 3.4
$ws.Range("Q36").Value = "28/10/2023 14:21"
$ws.Range("R36").Value = 1.77
$ws.Range("S36").Value = "27/10/2023 17:43"
$ws.Range("T36").Value = 1.74
$ws.Range("U36").Value = "28/10/2023 14:21"
$ws.Range("V36").Value = "https://www.betexplorer.com/football/malta/premier-league/valletta-hamrun/vkWrB6tA/"

# Row 37
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "malta"
$ws.Range("C37").Value = "premier-league"
$ws.Range("D37").Value = "2023-2024"
$ws.Range("E37").Value = 45227.72916666666
$ws.Range("F37").Value = "Naxxar"
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = "Hibernians"
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 6.73
$ws.Range("K37").Value = "27/10/2023 17:43"
$ws.Range("L37").Value = 5.51
$ws.Range("M37").Value = "28/10/2023 17:27"
$ws.Range("N37").Value = 4.85
$ws.Range("O37").Value = "27/10/2023 17:43"
$ws.Range("P37").Value = 4.36
$ws.Range("Q37").Value = "28/10/2023 17:27"
$ws.Range("R37").Value = 1.38
$ws.Range("S37").Value = "27/10/2023 17:43"
$ws.Range("T37").Value = 1.51
$ws.Range("U37").Value = "28/10/2023 17:27"
$ws.Range("V37").Value = "https://www.betexplorer.com/football/malta/premier-league/naxxar-lions-hibernians/0bXvCnR3/"

# Row 38
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "malta"
$ws.Range("C38").Value = "premier-league"
$ws.Range("D38").Value = "2023-2024"
$ws.Range("E38").Value = 45227.75
$ws.Range("F38").Value = "Mosta"
$ws.Range("G38").Value = 3
$ws.Range("H38").Value = "Marsaxlokk"
$ws.Range("I38").Value = 2
$ws.Range("J38").Value = 2.39
$ws.Range("K38").Value = "27/10/2023 05:12"
$ws.Range("L38").Value = 3.66
$ws.Range("M38").Value = "28/10/2023 17:55"
$ws.Range("N38").Value = 3.14
$ws.Range("O38").Value = "27/10/2023 05:12"
$ws.Range("P38").Value = 3.29
$ws.Range("Q38").Value = "28/10/2023 17:55"
$ws.Range("R38").Value = 2.72
$ws.Range("S38").Value = "27/10/2023 05:12"
$ws.Range("T38").Value = 2.02
$ws.Range("U38").Value = "28/10/2023 17:55"
$ws.Range("V38").Value = "https://www.betexplorer.com/football/malta/premier-league/mosta-fc-marsaxlokk/YBVnAQeG/"

# Row 39
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "malta"
$ws.Range("C39").Value = "premier-league"
$ws.Range("D39").Value = "2023-2024"
$ws.Range("E39").Value = 45228.45833333334
$ws.Range("F39").Value = "Sirens"
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = "Santa Lucia"
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2.45
$ws.Range("K39").Value = "28/10/2023 09:43"
$ws.Range("L39").Value = 2.51
$ws.Range("M39").Value = "29/10/2023 10:55"
$ws.Range("N39").Value = 3.08
$ws.Range("O39").Value = "28/10/2023 09:43"
$ws.Range("P39").Value = 2.86
$ws.Range("Q39").Value = "29/10/2023 10:55"
$ws.Range("R39").Value = 2.62
$ws.Range("S39").Value = "28/10/2023 09:43"
$ws.Range("T39").Value = 3.09
$ws.Range("U39").Value = "29/10/2023 10:55"
$ws.Range("V39").Value = "https://www.betexplorer.com/football/malta/premier-league/sirens-santa-lucia/nFZj9pBM/"

# Row 40
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "malta"
$ws.Range("C40").Value = "premier-league"
$ws.Range("D40").Value = "2023-2024"
$ws.Range("E40").Value = 45228.625
$ws.Range("F40").Value = "Balzan"
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = "Birkirkara"
$ws.Range("I40").Value = 3
$ws.Range("J40").Value = 2.8
$ws.Range("K40").Value = "28/10/2023 09:43"
$ws.Range("L40").Value = 3.6
$ws.Range("M40").Value = "29/10/2023 14:51"
$ws.Range("N40").Value = 3.03
$ws.Range("O40").Value = "28/10/2023 09:43"
$ws.Range("P40").Value = 2.56
$ws.Range("Q40").Value = "29/10/2023 14:51"
$ws.Range("R40").Value = 2.34
$ws.Range("S40").Value = "28/10/2023 09:43"
$ws.Range("T40").Value = 2.49
$ws.Range("U40").Value = "29/10/2023 14:51"
$ws.Range("V40").Value = "https://www.betexplorer.com/football/malta/premier-league/balzan-fc-birkirkara/WnBN3OYq/"

# Row 41
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "malta"
$ws.Range("C41").Value = "premier-league"
$ws.Range("D41").Value = "2023-2024"
$ws.Range("E41").Value = 45228.625
$ws.Range("F41").Value = "Gudja"
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = "Floriana"
$ws.Range("I41").Value = 1
$ws.Range("J41").Value = 4.82
$ws.Range("K41").Value = "28/10/2023 09:43"
$ws.Range("L41").Value = 9.96
$ws.Range("M41").Value = "29/10/2023 14:43"
$ws.Range("N41").Value = 3.62
$ws.Range("O41").Value = "28/10/2023 09:43"
$ws.Range("P41").Value = 5.01
$ws.Range("Q41").Value = "29/10/2023 14:43"
$ws.Range("R41").Value = 1.58
$ws.Range("S41").Value = "28/10/2023 09:43"
$ws.Range("T41").Value = 1.3
$ws.Range("U41").Value = "29/10/2023 10:32"
$ws.Range("V41").Value = "https://www.betexplorer.com/football/malta/premier-league/gudja-floriana/UwPe84QS/"

# Row 42
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "malta"
$ws.Range("C42").Value = "premier-league"
$ws.Range("D42").Value = "2023-2024"
$ws.Range("E42").Value = 45228.70833333334
$ws.Range("F42").Value = "Sliema"
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = "Gzira"
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 2.32
$ws.Range("K42").Value = "29/10/2023 13:42"
$ws.Range("L42").Value = 2.64
$ws.Range("M42").Value = "29/10/2023 15:00"
$ws.Range("N42").Value = 3.35
$ws.Range("O42").Value = "29/10/2023 13:42"
$ws.Range("P42").Value = 3.26
$ws.Range("Q42").Value = "29/10/2023 15:03"
$ws.Range("R42").Value = 2.79
$ws.Range("S42").Value = "29/10/2023 13:42"
$ws.Range("T42").Value = 2.56
$ws.Range("U42").Value = "29/10/2023 15:00"
$ws.Range("V42").Value = "https://www.betexplorer.com/football/malta/premier-league/sliema-gzira/0MFR2rlj/"

